$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values for the symbol list refresh.
# Force each target cell to Text format before assignment so that numeric-looking
# strings (e.g. "303.19") and percentages (e.g. "3.00%") are preserved verbatim
# as text, matching the original inline-string cell contents, then clear the
# temporary number-format override so no stray cell styling is introduced.
$updates = @(
    @{ Cell = 'D2'; Value = '303.19' }
    @{ Cell = 'E3'; Value = '3.00%' }
    @{ Cell = 'E4'; Value = '-3.67%' }
    @{ Cell = 'D5'; Value = '0.07833' }
    @{ Cell = 'E5'; Value = '-1.39%' }
    @{ Cell = 'D6'; Value = '2.036' }
    @{ Cell = 'E6'; Value = '-8.92%' }
    @{ Cell = 'D7'; Value = '7.834' }
    @{ Cell = 'E7'; Value = '0.72%' }
    @{ Cell = 'D8'; Value = '3.827' }
    @{ Cell = 'E8'; Value = '-0.85%' }
    @{ Cell = 'D9'; Value = '0.9218' }
    @{ Cell = 'E9'; Value = '-0.49%' }
    @{ Cell = 'E10'; Value = '1.77%' }
    @{ Cell = 'D11'; Value = '0.07889' }
    @{ Cell = 'E11'; Value = '5.84%' }
    @{ Cell = 'D12'; Value = '0.08588' }
    @{ Cell = 'E12'; Value = '-9.07%' }
    @{ Cell = 'D13'; Value = '0.03163' }
    @{ Cell = 'E13'; Value = '4.34%' }
    @{ Cell = 'D14'; Value = '0.1004' }
    @{ Cell = 'E14'; Value = '0.04%' }
    @{ Cell = 'D15'; Value = '0.001511' }
    @{ Cell = 'E15'; Value = '-0.30%' }
    @{ Cell = 'D16'; Value = '0.005901' }
    @{ Cell = 'E16'; Value = '0.73%' }
    @{ Cell = 'E17'; Value = '2,111.55%' }
    @{ Cell = 'D18'; Value = '3.465' }
    @{ Cell = 'E18'; Value = '-0.44%' }
    @{ Cell = 'E19'; Value = '-4.92%' }
    @{ Cell = 'D21'; Value = '0.1317' }
    @{ Cell = 'E21'; Value = '-2.23%' }
    @{ Cell = 'D22'; Value = '4.279' }
    @{ Cell = 'E22'; Value = '9.13%' }
    @{ Cell = 'D23'; Value = '0.1999' }
    @{ Cell = 'E23'; Value = '17.78%' }
    @{ Cell = 'D24'; Value = '0.04577' }
    @{ Cell = 'E24'; Value = '-0.77%' }
    @{ Cell = 'D25'; Value = '0.001223' }
    @{ Cell = 'E25'; Value = '-1.64%' }
    @{ Cell = 'D26'; Value = '0.004449' }
    @{ Cell = 'E26'; Value = '-0.60%' }
    @{ Cell = 'E27'; Value = '4.33%' }
    @{ Cell = 'D39'; Value = '0.01742' }
    @{ Cell = 'E39'; Value = '-1.29%' }
    @{ Cell = 'E40'; Value = '4.09%' }
    @{ Cell = 'D41'; Value = '0.007496' }
    @{ Cell = 'E41'; Value = '7.47%' }
    @{ Cell = 'D42'; Value = '0.1365' }
    @{ Cell = 'E42'; Value = '0.19%' }
    @{ Cell = 'E43'; Value = '7.93%' }
    @{ Cell = 'D44'; Value = '0.01056' }
    @{ Cell = 'E44'; Value = '10.38%' }
    @{ Cell = 'D45'; Value = '0.00006317' }
    @{ Cell = 'E45'; Value = '0.54%' }
    @{ Cell = 'E46'; Value = '0.15%' }
    @{ Cell = 'E47'; Value = '-61.07%' }
    @{ Cell = 'D48'; Value = '0.8205' }
    @{ Cell = 'E48'; Value = '9.91%' }
    @{ Cell = 'D49'; Value = '0.00002099' }
    @{ Cell = 'E49'; Value = '0.15%' }
    @{ Cell = 'E50'; Value = '0.15%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}

Write-Output ("Updated " + $updates.Count + " cells")
